$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing counts (QAPF map area 2a+2b)
$ws.Range("B2").Value = 149437
$ws.Range("B3").Value = 131

# Remove the last row (quartz monzodiorite / quartz monzogabbro) entirely
$ws.Range("A4:B4").ClearContents()
$ws.Rows.Item(4).Delete()
